$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7571855626538309
$ws.Range("C2").Value = 0.03726687579720078
$ws.Range("D2").Value = 0.1309344824700851
$ws.Range("F2").Value = 2.506004254945339
$ws.Range("G2").Value = 0.002551835714879697
$ws.Range("I2").Value = 2.074096297897299
$ws.Range("K2").Value = 0.5616931974663544
$ws.Range("L2").Value = 0.3067429217273911
$ws.Range("B3").Value = 0.7295575982021205
$ws.Range("C3").Value = 0.03284928380706731
$ws.Range("D3").Value = 0.1302543244145014
$ws.Range("F3").Value = 2.451772369326889
$ws.Range("G3").Value = 0.002556147667515718
$ws.Range("I3").Value = 2.039805160593204
$ws.Range("K3").Value = 0.5306056821829372
$ws.Range("L3").Value = 0.2961505466043661
$ws.Range("B4").Value = 0.7131966578589868
$ws.Range("C4").Value = 0.03012514075778938
$ws.Range("D4").Value = 0.1298242402744876
$ws.Range("F4").Value = 2.419280197422609
$ws.Range("G4").Value = 0.002558934494636622
$ws.Range("I4").Value = 2.019294305946659
$ws.Range("K4").Value = 0.5119651676941714
$ws.Range("L4").Value = 0.2898317989800319
$ws.Range("B5").Value = 0.7066809380389429
$ws.Range("C5").Value = 0.02901199767622131
$ws.Range("D5").Value = 0.1296458163119105
$ws.Range("F5").Value = 2.4062413926641
$ws.Range("G5").Value = 0.002560105285129671
$ws.Range("I5").Value = 2.011072234402135
$ws.Range("K5").Value = 0.5044812670954002
$ws.Range("L5").Value = 0.2873032706030045
$ws.Range("B6").Value = 0.7056081559613006
$ws.Range("C6").Value = 0.02882697548768931
$ws.Range("D6").Value = 0.1296159974296636
$ws.Range("F6").Value = 2.404088488355413
$ws.Range("G6").Value = 0.002560301819378515
$ws.Range("I6").Value = 2.009715180364438
$ws.Range("K6").Value = 0.5032453448828846
$ws.Range("L6").Value = 0.2868862104209313
$ws.Range("B7").Value = 0.7131081712661569
$ws.Range("C7").Value = 0.03011014091079289
$ws.Range("D7").Value = 0.1298218468220611
$ws.Range("F7").Value = 2.419103534597184
$ws.Range("G7").Value = 0.002558950141855792
$ws.Range("I7").Value = 2.019182869190189
$ws.Range("K7").Value = 0.5118637829522186
$ws.Range("L7").Value = 0.2897975105822184
$ws.Range("B8").Value = 0.7475342865280936
$ws.Range("C8").Value = 0.03574607107786676
$ws.Range("D8").Value = 0.1307025346961765
$ws.Range("F8").Value = 2.487137247429985
$ws.Range("G8").Value = 0.002553293640214947
$ws.Range("I8").Value = 2.06215946255719
$ws.Range("K8").Value = 0.5508812244738976
$ws.Range("L8").Value = 0.303052200111253
$ws.Range("B9").Value = 0.819834711385738
$ws.Range("C9").Value = 0.0467093999661472
$ws.Range("D9").Value = 0.1323319017321438
$ws.Range("F9").Value = 2.626992621845403
$ws.Range("G9").Value = 0.002543301044789079
$ws.Range("I9").Value = 2.150782656054901
$ws.Range("K9").Value = 0.6309593420217823
$ws.Range("L9").Value = 0.3305193967052986
$ws.Range("B10").Value = 0.8758939321584762
$ws.Range("C10").Value = 0.05471674873302845
$ws.Range("D10").Value = 0.1334712140938024
$ws.Range("F10").Value = 2.733742741454421
$ws.Range("G10").Value = 0.002536622565502184
$ws.Range("I10").Value = 2.218592216791009
$ws.Range("K10").Value = 0.6919951778439213
$ws.Range("L10").Value = 0.3516106244862272
$ws.Range("B11").Value = 0.9020399474405565
$ws.Range("C11").Value = 0.05835069097737744
$ws.Range("D11").Value = 0.1339773451447357
$ws.Range("F11").Value = 2.783190534666716
$ws.Range("G11").Value = 0.00253372675992439
$ws.Range("I11").Value = 2.250037349478575
$ws.Range("K11").Value = 0.7202470601459368
$ws.Range("L11").Value = 0.36140617841302
$ws.Range("B12").Value = 0.912033677128079
$ws.Range("C12").Value = 0.05972563325198621
$ws.Range("D12").Value = 0.1341672856529215
$ws.Range("F12").Value = 2.802043679710664
$ws.Range("G12").Value = 0.002532650532246517
$ws.Range("I12").Value = 2.262031549461767
$ws.Range("K12").Value = 0.7310156676190331
$ws.Range("L12").Value = 0.365144584592457
$ws.Range("B13").Value = 0.9098772189709337
$ws.Range("C13").Value = 0.05942956466491012
$ws.Range("D13").Value = 0.1341264546572454
$ws.Range("F13").Value = 2.797977596199672
$ws.Range("G13").Value = 0.002532881413709892
$ws.Range("I13").Value = 2.259444522474581
$ws.Range("K13").Value = 0.7286933269421922
$ws.Range("L13").Value = 0.3643381581342169
$ws.Range("B14").Value = 0.9028602768982239
$ws.Range("C14").Value = 0.05846383089397023
$ws.Range("D14").Value = 0.133993005991961
$ws.Range("F14").Value = 2.784739017304048
$ws.Range("G14").Value = 0.002533637810803241
$ws.Range("I14").Value = 2.251022380350264
$ws.Range("K14").Value = 0.7211315906837399
$ws.Range("L14").Value = 0.3617131565195564
$ws.Range("B15").Value = 0.8985742869789419
$ws.Range("C15").Value = 0.05787214362341331
$ws.Range("D15").Value = 0.1339110415758498
$ws.Range("F15").Value = 2.776646744107438
$ws.Range("G15").Value = 0.002534103773058519
$ws.Range("I15").Value = 2.245874871853545
$ws.Range("K15").Value = 0.716508964113757
$ws.Range("L15").Value = 0.3601090533833684
$ws.Range("B16").Value = 0.874198207009357
$ws.Range("C16").Value = 0.05447909583973853
$ws.Range("D16").Value = 0.1334378950631354
$ws.Range("F16").Value = 2.730529125013817
$ws.Range("G16").Value = 0.002536814666844065
$ws.Range("I16").Value = 2.216549291297682
$ws.Range("K16").Value = 0.6901586628053167
$ws.Range("L16").Value = 0.3509745188381146
$ws.Range("B17").Value = 0.8594094114892812
$ws.Range("C17").Value = 0.05239542140273556
$ws.Range("D17").Value = 0.1331445425481377
$ws.Range("F17").Value = 2.702465143992356
$ws.Range("G17").Value = 0.002538514074059796
$ws.Range("I17").Value = 2.198712633421366
$ws.Range("K17").Value = 0.674118383989935
$ws.Range("L17").Value = 0.3454223642748104
$ws.Range("B18").Value = 0.8509639320782298
$ws.Range("C18").Value = 0.05119612824202591
$ws.Range("D18").Value = 0.1329746689505491
$ws.Range("F18").Value = 2.686406887295362
$ws.Range("G18").Value = 0.002539504925489444
$ws.Range("I18").Value = 2.188509731930466
$ws.Range("K18").Value = 0.6649382016664731
$ws.Range("L18").Value = 0.3422478274350027
$ws.Range("B19").Value = 0.8481148495863806
$ws.Range("C19").Value = 0.05078992484565958
$ws.Range("D19").Value = 0.1329169552065075
$ws.Range("F19").Value = 2.680984136437871
$ws.Range("G19").Value = 0.002539842714951339
$ws.Range("I19").Value = 2.185064851788866
$ws.Range("K19").Value = 0.661837800530634
$ws.Range("L19").Value = 0.3411762286004318
$ws.Range("B20").Value = 0.8609774265695762
$ws.Range("C20").Value = 0.05261731619015109
$ws.Range("D20").Value = 0.1331758887340477
$ws.Range("F20").Value = 2.705443964738237
$ws.Range("G20").Value = 0.002538331783391798
$ws.Range("I20").Value = 2.200605547907983
$ws.Range("K20").Value = 0.6758211607914006
$ws.Range("L20").Value = 0.3460114419612523
$ws.Range("B21").Value = 0.9049188031338531
$ws.Range("C21").Value = 0.0587475208783701
$ws.Range("D21").Value = 0.1340322495784889
$ws.Range("F21").Value = 2.788624019003493
$ws.Range("G21").Value = 0.002533415087055524
$ws.Range("I21").Value = 2.253493812364724
$ws.Range("K21").Value = 0.7233507474454939
$ws.Range("L21").Value = 0.3624833938561096
$ws.Range("B22").Value = 0.9341780199895311
$ws.Range("C22").Value = 0.06274728691796838
$ws.Range("D22").Value = 0.1345819156711272
$ws.Range("F22").Value = 2.843735445989267
$ws.Range("G22").Value = 0.002530320312794519
$ws.Range("I22").Value = 2.288564458612001
$ws.Range("K22").Value = 0.7548235928182692
$ws.Range("L22").Value = 0.3734181114922421
$ws.Range("B23").Value = 0.9185122767970313
$ws.Range("C23").Value = 0.06061311759272314
$ws.Range("D23").Value = 0.1342894562005057
$ws.Range("F23").Value = 2.81425267661038
$ws.Range("G23").Value = 0.002531961238249509
$ws.Range("I23").Value = 2.269800182512171
$ws.Range("K23").Value = 0.7379883733370605
$ws.Range("L23").Value = 0.3675665076043799
$ws.Range("B24").Value = 0.8602683500510295
$ws.Range("C24").Value = 0.05251700180376417
$ws.Range("D24").Value = 0.1331617209294862
$ws.Range("F24").Value = 2.704097002893803
$ws.Range("G24").Value = 0.002538414153912577
$ws.Range("I24").Value = 2.199749600482861
$ws.Range("K24").Value = 0.6750512058840172
$ws.Range("L24").Value = 0.3457450655328671
$ws.Range("B25").Value = 0.7997606573858036
$ws.Range("C25").Value = 0.04375233937963685
$ws.Range("D25").Value = 0.1319013882480782
$ws.Range("F25").Value = 2.588461065421257
$ws.Range("G25").Value = 0.002545887331755091
$ws.Range("I25").Value = 2.126337435188816
$ws.Range("K25").Value = 0.6089114528705863
$ws.Range("L25").Value = 0.3229296730726503
